# Commit: "add retA_real to cpars"
# Adds 6 new rows (L5_y, LFS_y, Vmaxlen_y, LR5_y, LFR_y, Rmaxlen_y) to the
# "cpars" sheet, describing the "by sim and year" variants of the
# selectivity/retention parameters, following the same Slot/Class/Desc/Type
# layout used by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cpars")

# Column B (Class) is the same "numeric matrix" type for every new row, and
# column D (Type) is "Fleet" for all of them - same values already used
# elsewhere on this sheet (e.g. row 2's "Karray"-style entries use the same
# matrix-dim description).
$matrixClass = "numeric matrix dim=c(nsim, nyears+proyears)"
$fleetType   = "Fleet"

$slots = @("L5_y", "LFS_y", "Vmaxlen_y", "LR5_y", "LFR_y", "Rmaxlen_y")
$descs = @(
  "length-at-5 percent selection by sim and year",
  "length at full selection by sim and year",
  "vulnerability at max length by sim and year",
  "length-at-5 percent retention by sim and year",
  "length at full retention by sim and year",
  "retention at max length by sim and year"
)

$startRow = 103

# Fill column by column (A, then B, then C, then D) so new shared strings are
# interned in the same order the original workbook used.
for ($i = 0; $i -lt $slots.Length; $i++) {
  $ws.Cells.Item($startRow + $i, 1).Value = $slots[$i]
}
for ($i = 0; $i -lt $slots.Length; $i++) {
  $ws.Cells.Item($startRow + $i, 2).Value = $matrixClass
}
for ($i = 0; $i -lt $slots.Length; $i++) {
  $ws.Cells.Item($startRow + $i, 3).Value = $descs[$i]
}
for ($i = 0; $i -lt $slots.Length; $i++) {
  $ws.Cells.Item($startRow + $i, 4).Value = $fleetType
}

# Match the saved view state: selection moved to C109 after the new rows.
$ws.Select()
$ws.Range("C109").Select()
